$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 33534314
$ws.Range("I70").Value = 41917590
$ws.Range("J70").Value = 1200
$ws.Range("K70").Value = 125752770
$ws.Range("L70").Value = 3600
$ws.Range("M70").Value = -125752500
$ws.Range("N70").Value = -4140

$ws.Range("H73").Value = 33534314
$ws.Range("I73").Value = 41917590
$ws.Range("J73").Value = 1200
$ws.Range("K73").Value = 125752770
$ws.Range("L73").Value = 3600
$ws.Range("M73").Value = -125751834
$ws.Range("N73").Value = -5472

$ws.Range("H86").Value = 1279
$ws.Range("J86").Value = 1497.5
$ws.Range("L86").Value = 1497.5
$ws.Range("N86").Value = -3743.5

$ws.Range("H89").Value = 1279
$ws.Range("J89").Value = 1497.5
$ws.Range("L89").Value = 7487.5
$ws.Range("N89").Value = -18719.5

$ws.Range("H112").Value = 2748.5334
$ws.Range("J112").Value = 2808.8276
$ws.Range("L112").Value = 8426.4828
$ws.Range("N112").Value = -10642.4828

$ws.Range("H137").Value = 2017.0577
$ws.Range("I137").Value = 1585.5454
$ws.Range("K137").Value = 4756.6362
$ws.Range("M137").Value = -2206.6362

$ws.Range("H138").Value = 3273.63
$ws.Range("I138").Value = 2007.3478
$ws.Range("J138").Value = 3651.87
$ws.Range("K138").Value = 6022.0434
$ws.Range("L138").Value = 10955.61
$ws.Range("M138").Value = -882.0434000000005
$ws.Range("N138").Value = -21235.61

$ws.Range("H141").Value = 5239.7144
$ws.Range("I141").Value = 1730.0741
$ws.Range("K141").Value = 5190.2223
$ws.Range("M141").Value = -10.22230000000036

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10637.878
$ws.Range("I32").Value = 11544.339
$ws.Range("J32").Value = 5954.5
$ws.Range("K32").Value = 11544.339
$ws.Range("L32").Value = 5954.5
$ws.Range("M32").Value = -11257.339
$ws.Range("N32").Value = -6528.5

$ws.Range("H74").Value = 1454.6562
$ws.Range("I74").Value = 1255.3214
$ws.Range("J74").Value = 2850
$ws.Range("K74").Value = 1255.3214
$ws.Range("L74").Value = 2850
$ws.Range("M74").Value = -381.3214
$ws.Range("N74").Value = -4598

$ws.Range("H77").Value = 1454.6562
$ws.Range("I77").Value = 1255.3214
$ws.Range("J77").Value = 2850
$ws.Range("K77").Value = 6276.607
$ws.Range("L77").Value = 14250
$ws.Range("M77").Value = -1908.607
$ws.Range("N77").Value = -22986

$ws.Range("H122").Value = 5994.6665
$ws.Range("I122").Value = 6988
$ws.Range("J122").Value = 2518
$ws.Range("K122").Value = 20964
$ws.Range("L122").Value = 7554
$ws.Range("M122").Value = -18514
$ws.Range("N122").Value = -12454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 33804.28
$ws.Range("I86").Value = 2364.3157
$ws.Range("J86").Value = 79755
$ws.Range("K86").Value = 2364.3157
$ws.Range("L86").Value = 79755
$ws.Range("M86").Value = -1241.3157
$ws.Range("N86").Value = -82001

$ws.Range("H89").Value = 33804.28
$ws.Range("I89").Value = 2364.3157
$ws.Range("J89").Value = 79755
$ws.Range("K89").Value = 11821.5785
$ws.Range("L89").Value = 398775
$ws.Range("M89").Value = -6205.5785
$ws.Range("N89").Value = -410007

$ws.Range("H94").Value = 886.5
$ws.Range("I94").Value = 759.8182
$ws.Range("J94").Value = 1085.5714
$ws.Range("K94").Value = 759.8182
$ws.Range("L94").Value = 1085.5714
$ws.Range("M94").Value = -308.8182
$ws.Range("N94").Value = -1987.5714

$ws.Range("H134").Value = 2503.889
$ws.Range("I134").Value = 1942.7273
$ws.Range("J134").Value = 3385.7144
$ws.Range("K134").Value = 5828.1819
$ws.Range("L134").Value = 10157.1432
$ws.Range("M134").Value = -3293.1819
$ws.Range("N134").Value = -15227.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4155.852
$ws.Range("I31").Value = 3068.842
$ws.Range("J31").Value = 6737.5
$ws.Range("K31").Value = 3068.842
$ws.Range("L31").Value = 6737.5
$ws.Range("M31").Value = -2773.842
$ws.Range("N31").Value = -7327.5

$ws.Range("H34").Value = 4155.852
$ws.Range("I34").Value = 3068.842
$ws.Range("J34").Value = 6737.5
$ws.Range("K34").Value = 3068.842
$ws.Range("L34").Value = 6737.5
$ws.Range("M34").Value = -2866.842
$ws.Range("N34").Value = -7141.5

$ws.Range("H58").Value = 4782.4
$ws.Range("I58").Value = 7506
$ws.Range("J58").Value = 2966.6667
$ws.Range("K58").Value = 7506
$ws.Range("L58").Value = 2966.6667
$ws.Range("M58").Value = -7303
$ws.Range("N58").Value = -3372.6667

$ws.Range("H99").Value = 1629
$ws.Range("I99").Value = 1677.2858
$ws.Range("J99").Value = 1493.8
$ws.Range("K99").Value = 1677.2858
$ws.Range("L99").Value = 1493.8
$ws.Range("M99").Value = -179.2858000000001
$ws.Range("N99").Value = -4489.8

$ws.Range("H126").Value = 1629
$ws.Range("I126").Value = 1677.2858
$ws.Range("J126").Value = 1493.8
$ws.Range("K126").Value = 5031.857400000001
$ws.Range("L126").Value = 4481.4
$ws.Range("M126").Value = -2561.857400000001
$ws.Range("N126").Value = -9421.4

$ws.Range("H131").Value = 24250.5
$ws.Range("J131").Value = 24250.5
$ws.Range("L131").Value = 24250.5
$ws.Range("N131").Value = -34330.5

$ws.Range("H134").Value = 2529.6875
$ws.Range("I134").Value = 1343
$ws.Range("K134").Value = 4029
$ws.Range("M134").Value = -1494

$ws.Range("H136").Value = 4782.4
$ws.Range("I136").Value = 7506
$ws.Range("J136").Value = 2966.6667
$ws.Range("K136").Value = 22518
$ws.Range("L136").Value = 8900.000100000001
$ws.Range("M136").Value = -19968
$ws.Range("N136").Value = -14000.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1122.931
$ws.Range("I5").Value = 1535.2778
$ws.Range("J5").Value = 448.18182
$ws.Range("K5").Value = 4605.8334
$ws.Range("L5").Value = 1344.54546
$ws.Range("M5").Value = -4493.8334
$ws.Range("N5").Value = -1568.54546

$ws.Range("H26").Value = 565
$ws.Range("J26").Value = 866.6667
$ws.Range("L26").Value = 2600.0001
$ws.Range("N26").Value = -3176.0001

$ws.Range("H131").Value = 19611000
$ws.Range("I131").Value = 815.8
$ws.Range("K131").Value = 2447.4
$ws.Range("M131").Value = 2592.6

$ws.Range("H132").Value = 1644.3334
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 1662.375
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 14961.375
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -20021.375

$ws.Range("H134").Value = 3950.4736
$ws.Range("I134").Value = 2045.619
$ws.Range("J134").Value = 6303.5293
$ws.Range("K134").Value = 6136.857
$ws.Range("L134").Value = 18910.5879
$ws.Range("M134").Value = -1066.857
$ws.Range("N134").Value = -29050.5879

$ws.Range("H135").Value = 1122.931
$ws.Range("I135").Value = 1535.2778
$ws.Range("J135").Value = 448.18182
$ws.Range("K135").Value = 13817.5002
$ws.Range("L135").Value = 4033.63638
$ws.Range("M135").Value = -11282.5002
$ws.Range("N135").Value = -9103.63638

$ws.Range("H139").Value = 1919.1666
$ws.Range("I139").Value = 1453.75
$ws.Range("J139").Value = 2850
$ws.Range("K139").Value = 4361.25
$ws.Range("L139").Value = 8550
$ws.Range("M139").Value = 778.75
$ws.Range("N139").Value = -18830

$ws.Range("H140").Value = 5212.353
$ws.Range("I140").Value = 1360
$ws.Range("J140").Value = 7909
$ws.Range("K140").Value = 4080
$ws.Range("L140").Value = 23727
$ws.Range("M140").Value = 1100
$ws.Range("N140").Value = -34087

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2375.85
$ws.Range("I126").Value = 1792.2727
$ws.Range("J126").Value = 3089.111
$ws.Range("K126").Value = 5376.8181
$ws.Range("L126").Value = 9267.332999999999
$ws.Range("M126").Value = -2906.8181
$ws.Range("N126").Value = -14207.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3543.923
$ws.Range("I16").Value = 674.2
$ws.Range("J16").Value = 5337.5
$ws.Range("K16").Value = 674.2
$ws.Range("L16").Value = 5337.5
$ws.Range("M16").Value = -504.2
$ws.Range("N16").Value = -5677.5

$ws.Range("H122").Value = 150002000
$ws.Range("I122").Value = 125002990
$ws.Range("J122").Value = 200000000
$ws.Range("K122").Value = 375008970
$ws.Range("L122").Value = 600000000
$ws.Range("M122").Value = -375006520
$ws.Range("N122").Value = -600004900

$ws.Range("H132").Value = 4720.3237
$ws.Range("I132").Value = 4936.222
$ws.Range("J132").Value = 3887.5715
$ws.Range("K132").Value = 14808.666
$ws.Range("L132").Value = 11662.7145
$ws.Range("M132").Value = -12278.666
$ws.Range("N132").Value = -16722.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 35966.668
$ws.Range("I54").Value = 30000
$ws.Range("J54").Value = 38950
$ws.Range("K54").Value = 30000
$ws.Range("L54").Value = 38950
$ws.Range("M54").Value = -29480
$ws.Range("N54").Value = -39990

$ws.Range("H96").Value = 1166.6666
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H136").Value = 2385.3333
$ws.Range("I136").Value = 2703.6428
$ws.Range("J136").Value = 2106.8125
$ws.Range("K136").Value = 8110.928400000001
$ws.Range("L136").Value = 6320.4375
$ws.Range("M136").Value = -5560.928400000001
$ws.Range("N136").Value = -11420.4375
